# updated the sheetname as per feedback
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (Sheet1 -> DemonstrationTest)
$ws.Name = "DemonstrationTest"

# Excel constant values (numeric, since named constants aren't resolved by this host)
$xlLeft = -4131
$xlCenter = -4108

# --- Header row (row 1) ---
# A1 & B1 keep their existing bold/centered/bordered look, but the font
# family moves from Cambria to Arial and the cells gain wrap text.
$hdrRange = $ws.Range("A1:B1")
$hdrRange.Font.Name = "Arial"
$hdrRange.WrapText = $true

# B1's font color becomes an explicit black (was theme color 1).
$ws.Range("B1").Font.Color = 0

# C1 & D1 go from an unformatted/unbordered cell to match the bold Arial
# header look, but left-aligned with a thin border and wrap text.
$cdRange = $ws.Range("C1:D1")
$cdRange.Font.Name = "Arial"
$cdRange.Font.Size = 8
$cdRange.Font.Bold = $true
$cdRange.Borders.LineStyle = "xlContinuous"
$cdRange.Borders.Weight = "xlThin"
$cdRange.Borders.Color = 0
$cdRange.HorizontalAlignment = $xlLeft
$cdRange.WrapText = $true

# E1 goes from a plain (non-bold, theme-colored, unbordered) cell to the
# same bold/black/bordered/left-aligned look as B1, but left-aligned.
$e1 = $ws.Range("E1")
$e1.Font.Name = "Arial"
$e1.Font.Size = 8
$e1.Font.Bold = $true
$e1.Font.Color = 0
$e1.Borders.LineStyle = "xlContinuous"
$e1.Borders.Weight = "xlThin"
$e1.Borders.Color = 0
$e1.HorizontalAlignment = $xlLeft
$e1.WrapText = $true

# --- Data row (row 2) ---
# A2 & B2 keep their existing centered/bordered look, just gain wrap text.
$ws.Range("A2:B2").WrapText = $true

# C2, D2 & E2 go from unformatted/unbordered cells to a bordered,
# left-aligned, wrap-text look (same font as A2/B2 already have).
$bodyRange = $ws.Range("C2:E2")
$bodyRange.Borders.LineStyle = "xlContinuous"
$bodyRange.Borders.Weight = "xlThin"
$bodyRange.Borders.Color = 0
$bodyRange.HorizontalAlignment = $xlLeft
$bodyRange.WrapText = $true

Write-Host "done"
